$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a data row (Conta as text to preserve leading zeros, Nome as text, Saldo as number)
function Set-DataRow($row, $conta, $nome, $saldo) {
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $conta
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $nome
    $ws.Cells.Item($row, 3).Value = $saldo
}

# 1) Insert new row for account 005696533 / CARLOS / 50000 right after BLUEMETRIX (row 3),
#    i.e. immediately before RENATO (originally row 4).
$ws.Rows.Item(4).Insert()
Set-DataRow 4 "005696533" "CARLOS" 50000

# 2) Move the block RAFAEL(004222784)/THIAGO/SOLANGE/ISABEL that currently sits right after
#    EDWARD (row 7) to instead sit right BEFORE EDWARD. Delete the existing 4 rows (8-11,
#    highest index first) and re-insert them as 4 fresh rows starting at row 7.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

Set-DataRow 7 "004222784" "RAFAEL" 24105.45
Set-DataRow 8 "005064129" "THIAGO" 23901.33
Set-DataRow 9 "004455489" "SOLANGE" 21104.41
Set-DataRow 10 "005624730" "ISABEL" 20000

# 3) EDWARD (now at row 11) balance drops from 30814.68 to 18000.
$ws.Cells.Item(11, 3).Value = 18000

# 4) Remove ADELE (004575632) entirely - row 14.
$ws.Rows.Item(14).Delete()

# 5) Insert new row for account 005529100 / DIMITRI / 8000 right after DANIELA (now row 14).
$ws.Rows.Item(15).Insert()
Set-DataRow 15 "005529100" "DIMITRI" 8000
